# Update Financials figures for MDB Yearly sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 166000
$ws.Range("D10").Value = 123200
$ws.Range("D17").Value = 250900
$ws.Range("D18").Value = -84900
$ws.Range("D21").Value = -79000
$ws.Range("D23").Value = -82700
$ws.Range("D26").Value = -84000
$ws.Range("D27").Value = -84000
$ws.Range("D33").Value = -84000
$ws.Range("D35").Value = -84000
$ws.Range("D45").Value = 21400
$ws.Range("D46").Value = 347300
$ws.Range("D48").Value = 119100
$ws.Range("D49").Value = 5000
$ws.Range("D52").Value = 31100
$ws.Range("D54").Value = 432800
$ws.Range("D59").Value = 110300
$ws.Range("D60").Value = 112500
$ws.Range("D62").Value = 72700
$ws.Range("D66").Value = 185200
$ws.Range("D72").Value = -389600
$ws.Range("D76").Value = 247700
$ws.Range("D81").Value = -84000
